# TIMES_Energy_HROR.xlsx - update EV profile source data
#
# The sheet used to report PJ figures for the year "2050"; it is now
# re-labelled as an activity-unit table ("Unit - activity1" / "PJ")
# sourced from the new RAMP-mobility results, per the commit message
# ("Updated EV profiles source with RAMP-mobility results").
#
#   B1: "PJ"   -> "Unit - activity1"   (column/table caption)
#   B2: "2050" -> "PJ"                 (unit row, now reads "PJ")
#   A2:A25 (Region / country codes) and the B3:B25 numeric data are
#   unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the header / unit labels -------------------------------------
# Leading apostrophes preserve the "quote prefix" text formatting that was
# already applied to these label cells (so they keep being stored as
# left/''-prefixed text rather than being re-interpreted).
$ws.Range("B2").Value = "'PJ"
$ws.Range("B1").Value = "'Unit - activity1"

# Give the new table caption ("Unit - activity1") its own distinct look:
# centered both horizontally and vertically, using the darker label font
# already used elsewhere in the sheet.
$ws.Range("B1").Font.Color = 3355443
$ws.Range("B1").VerticalAlignment = -4108

# --- Column B is narrower now that "2050" is no longer the widest entry --
$ws.Columns.Item(2).ColumnWidth = 10.7

# --- Clear the stale selection left over on cell C9 -----------------------
$null = $ws.Range("A1").Select()
